$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the x2 throughput/measurement error (doubling the values) ---
$ws.Range("F19").Value = 9.1866474
$ws.Range("H19").Value = 4.2191409999999996
$ws.Range("F20").Value = 11.4112834
$ws.Range("H20").Value = 9.4160400000000006
$ws.Range("F21").Value = 14.3144092
$ws.Range("H21").Value = 12.3920794
$ws.Range("F22").Value = 17.382556000000001
$ws.Range("H22").Value = 13.506182600000001
$ws.Range("F23").Value = 19.179984000000001
$ws.Range("H23").Value = 13.789209400000001

# --- Update the off-campus laptop owner's name ---
$ws.Range("B6").Value = "John Santos' Laptop (off campus)"

# --- Remove the stray "trial 3" label row (B40) entirely ---
$ws.Rows.Item(40).Delete()

# --- Remove the embedded picture/chart image ---
$ws.Shapes.Item(1).Delete()

# --- Update the active selection to reflect where the author ended up ---
$ws.Range("B40").Select()
